$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "VALOR MORA" total
$ws.Range("E11").Value = 148992

# 2. Update worker / period counts
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# 3. Update Salario Basico for existing row (period 2012) for worker JOHANNY
$ws.Range("G16").Value = 1423500

# 4. Insert two new data rows after row 16 (pushes the signature block down)
$ws.Rows("17:18").Insert()

# 5. Row 17: new worker GUILLERMO PANCRACIO ESPINOSA TORRES, period 2508
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "9293136"
$ws.Range("D17").Value = "GUILLERMO PANCRACIO ESPINOSA TORRES"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 644350

# 6. Row 18: existing worker JOHANNY MARIA NAVARRO SEGRERA, new period 2508
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "32936083"
$ws.Range("D18").Value = "JOHANNY MARIA NAVARRO SEGRERA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Column D is widest now to fit the new worker's longer name
$ws.Columns("D").ColumnWidth = 39.8
